# Append a new result row (row 15) to the "Sonuçlar" sheet, mirroring the
# structure of the existing rows (A:C text, D:M numbers, N:Q percent text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Cells.Item($row, 1).Value = "2025-08-12 08:50"
$ws.Cells.Item($row, 2).Value = "Deneme Testi Örnek"
$ws.Cells.Item($row, 3).Value = "Deneme Testi Örnek"
$ws.Cells.Item($row, 4).Value = 29
$ws.Cells.Item($row, 5).Value = 38
$ws.Cells.Item($row, 6).Value = 32
$ws.Cells.Item($row, 7).Value = 35
$ws.Cells.Item($row, 8).Value = 33.5
$ws.Cells.Item($row, 9).Value = 33.5
$ws.Cells.Item($row, 10).Value = 0.49
$ws.Cells.Item($row, 11).Value = 0.21
$ws.Cells.Item($row, 12).Value = 0.21
$ws.Cells.Item($row, 13).Value = 0.09
$ws.Cells.Item($row, 14).Value = "%48.71"
$ws.Cells.Item($row, 15).Value = "%21.08"
$ws.Cells.Item($row, 16).Value = "%21.08"
$ws.Cells.Item($row, 17).Value = "%9.13"
